$d = $word.ActiveDocument

# --- Bold the job-title portion of the three "Previous Positions" lines,
#     leaving the trailing " (date range)" text non-bold. Word's Find
#     naturally breaks the run into a bold run + a following plain run
#     when we set Bold only on the matched sub-range.

$titles = @(
    "Associate Principal Front-End Architect",
    "Lead Front-End Developer",
    "Senior Front-End Application Developer"
)

foreach ($title in $titles) {
    $rng = $d.Content
    $found = $rng.Find.Execute($title, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
    }
}

# --- Tighten the wording of the bullet about new feature development.

$d.Content.Find.Execute(
    "Developed new features for existing applications, updated existing functionalities, remediated defects, and rebuilt out-of-date assemblies and applications using more modern technologies",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Developed new features, updated existing features, remediated defects, and rebuilt out-of-date assemblies and applications using more modern technologies",
    2)
